$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-02-10 Monday" "2025-02-11 Tuesday"

Replace-Text "749×5=3745" "544×9=4896"
Replace-Text "316×4=1264" "159×2=318"
Replace-Text "340×8=2720" "664×8=5312"
Replace-Text "813×9=7317" "717×7=5019"
Replace-Text "357×4=1428" "211×3=633"

Replace-Text "345×3=1035" "383×2=766"
Replace-Text "863×9=7767" "185×5=925"
Replace-Text "264×2=528" "642×5=3210"
Replace-Text "858×8=6864" "605×6=3630"
Replace-Text "411×9=3699" "227×2=454"

Replace-Text "330×8=2640" "620×9=5580"
Replace-Text "176×4=704" "217×3=651"
Replace-Text "874×2=1748" "313×5=1565"
Replace-Text "956×6=5736" "327×7=2289"
Replace-Text "170×9=1530" "633×3=1899"

Replace-Text "102×8=816" "617×7=4319"
Replace-Text "870×3=2610" "508×4=2032"
Replace-Text "440×9=3960" "341×5=1705"
Replace-Text "417×3=1251" "315×3=945"
Replace-Text "113×9=1017" "793×2=1586"

Replace-Text "353×6=2118" "414×5=2070"
Replace-Text "819×7=5733" "324×5=1620"
Replace-Text "229×2=458" "434×7=3038"
Replace-Text "629×8=5032" "710×5=3550"
Replace-Text "186×2=372" "217×9=1953"
